{"js": "// Apply the commit's text-level edits using Word JavaScript API (Office.js)\n// search()+insertText(..., Word.InsertLocation.replace) on context.document.body.\n\nasync function replaceOnce(body, searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nasync function replaceAll(body, searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) ASUNTO row: \"...PRECIOS AL PEDIDO DE COMPRA N\u00ba 000022\" -> append new article.\nawait replaceOnce(\n  body,\n  \"PRECIOS AL PEDIDO DE COMPRA N\u00ba 000022\",\n  \"PRECIOS AL PEDIDO DE COMPRA N\u00ba 000078 - CEMENTO PORTLAND TIPO IP\"\n);\n\n// 2) REFERENCIA row (A): insert \"con \" before \"CUI: N\u00ba 2615887\" that follows\n//    the closing quote with a space (\" \u201d CUI: N\u00ba 2615887\") -- unique occurrence.\nawait replaceOnce(\n  body,\n  \"DEPARTAMENTO DE PUNO \u201d CUI: N\u00ba 2615887\",\n  \"DEPARTAMENTO DE PUNO \u201d con CUI: N\u00ba 2615887\"\n);\n\n// 3) 5) 8) All three identical mentions of the old INFORME number become the new one.\nawait replaceAll(\n  body,\n  \"INFORME N\u00ba 008-2025-MDSM/GDTI/SGI/SO-LEAB\",\n  \"INFORME N\u00ba 054-2026-MDSM/GDTI/SGI/RO-AICH\"\n);\n\n// 4) FECHA row date update.\nawait replaceOnce(\n  body,\n  \"San Miguel, 03 de Febrero de 2026\",\n  \"San Miguel, 12 de febrero de 2026\"\n);\n\n// 6) Remove duplicated \"de obra\".\nawait replaceOnce(\n  body,\n  \"con aprobaci\u00f3n del supervisor de obra de obra \",\n  \"con aprobaci\u00f3n del supervisor de obra \"\n);\n\n// 7) VALIDAN LA COTIZACI\u00d3N sentence: new purchase order number + trailing article.\nawait replaceOnce(\n  body,\n  \"COTIZACI\u00d3N AL PEDIDO DE COMPRA N\u00ba 000022 tambi\u00e9n\",\n  \"COTIZACI\u00d3N AL PEDIDO DE COMPRA N\u00ba 000078 DE CEMENTO PORTLAND TIPO IP tambi\u00e9n\"\n);\n", "ps1": "# Apply the commit's text-level edits using Word COM interop (PowerShell-style).\n# $word.ActiveDocument is pre-seeded by the harness.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text {\n    param(\n        [string]$Find,\n        [string]$ReplaceWith,\n        [int]$ReplaceScope = 2   # 2 = wdReplaceAll, 1 = wdReplaceOne\n    )\n    $rng = $d.Content.Duplicate\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($Find, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceWith, $ReplaceScope) | Out-Null\n}\n\n# 1) ASUNTO row: append the new purchase-order number and article description.\nReplace-Text \"PRECIOS AL PEDIDO DE COMPRA N\u00ba 000022\" \"PRECIOS AL PEDIDO DE COMPRA N\u00ba 000078 - CEMENTO PORTLAND TIPO IP\" 1\n\n# 2) REFERENCIA row (A): insert \"con \" before \"CUI: N\u00ba 2615887\" following the closing quote.\nReplace-Text \"DEPARTAMENTO DE PUNO \u201d CUI: N\u00ba 2615887\" \"DEPARTAMENTO DE PUNO \u201d con CUI: N\u00ba 2615887\" 1\n\n# 3) 5) 8) All three identical mentions of the old INFORME number become the new one.\nReplace-Text \"INFORME N\u00ba 008-2025-MDSM/GDTI/SGI/SO-LEAB\" \"INFORME N\u00ba 054-2026-MDSM/GDTI/SGI/RO-AICH\" 2\n\n# 4) FECHA row date update.\nReplace-Text \"San Miguel, 03 de Febrero de 2026\" \"San Miguel, 12 de febrero de 2026\" 1\n\n# 6) Remove duplicated \"de obra\".\nReplace-Text \"con aprobaci\u00f3n del supervisor de obra de obra \" \"con aprobaci\u00f3n del supervisor de obra \" 1\n\n# 7) VALIDAN LA COTIZACI\u00d3N sentence: new purchase-order number + trailing article.\nReplace-Text \"COTIZACI\u00d3N AL PEDIDO DE COMPRA N\u00ba 000022\" \"COTIZACI\u00d3N AL PEDIDO DE COMPRA N\u00ba 000078 DE CEMENTO PORTLAND TIPO IP\" 1\n"}
